
$d = $word.ActiveDocument

# --- Locate the last paragraph of the document (end of the 2024-05-17 entry) ---
$lastPara = $d.Paragraphs.Last

# --- Create the three new empty paragraphs first (blank separator, date
#     heading, body) so that later formatting (bold on the date line) does
#     not get inherited forward into the paragraphs created after it. ------
$lastPara.Range.InsertParagraphAfter()
$blankPara = $d.Paragraphs.Last
$blankPara.Range.InsertParagraphAfter()
$datePara = $d.Paragraphs.Last
$datePara.Range.InsertParagraphAfter()
$bodyPara = $d.Paragraphs.Last

# --- 1) Clean up the blank separator paragraph ------------------------------
# The engine seeds a brand-new empty paragraph with a placeholder run; typing
# a character then deleting just that character collapses the paragraph back
# down to a clean "no runs" state, matching how Word normally leaves an
# untouched blank paragraph.
$blankPara = $d.Paragraphs.Item(($d.Paragraphs.Count) - 2)
$blankPara.Range.InsertBefore("X")
$blankPara = $d.Paragraphs.Item(($d.Paragraphs.Count) - 2)
$placeholder = $d.Range($blankPara.Range.Start, $blankPara.Range.Start + 1)
$placeholder.Delete()

# --- 2) Fill in the bold date heading paragraph ------------------------------
$datePara = $d.Paragraphs.Item(($d.Paragraphs.Count) - 1)
$datePara.Range.InsertBefore('2024-05-21')
$datePara = $d.Paragraphs.Item(($d.Paragraphs.Count) - 1)
$datePara.Range.Font.Bold = $true
$datePara.Range.Font.BoldBi = $true

# --- 3) Fill in the journal body paragraph -----------------------------------
$bodyPara = $d.Paragraphs.Last
$bodyPara.Range.InsertBefore('First day on site.  This is exciting because coming into this job I didn''t really expect to be able to go to the sites of the companies I''d be working for, so this feels cool. It was also nice to meet the faces of those I had spent training sessions over teams with and really get a feel for the people on the project. It was pretty refreshing to understand how frustrated they are with their old historian system and really want something new, robust and effective. They expressed how they''ve been banging their heads against the wall and dying for this new system but there were a lot of factors preventing them from reaching this goal. It is a bit odd though because they did have this system in place originally but the way it was configured by a different integrator caused data retrieval to be slow and optimized or nonfunctional. This led them to just use their old system, but the issue is that the previous supporting company behind this old system no longer does business with them, so if anything were to happen to that system, they would be in the hole. Apparently, the company I''m working for never got wind of wanting to revamp their system so never got a bid for the work. Then we arrive in the modern day where we are now fixing the long trail of mistakes and it seems like those on the project are finally happy with the new system despite it not even being fully ready for use. Tomorrow is day 2 of the training so I''ll be back on site.')
